$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the shared "Massachusettes" string (affects A1:A4)
$ws.Range("A1:A4").Value = "Massachusetts {{ ma }}"

# Add a new numeric data row (row 10) - demonstrates support for
# non-string (numeric) cell types
$ws.Cells.Item(10, 1).Value = 2019
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(10, 3).Value = 16

# Update the active selection
$ws.Range("B15").Select()
